$d = $word.ActiveDocument

# Locate the "Status Cocktail" row in the first table and update the
# last cell's text from "error" to "ok".
$table = $d.Tables.Item(1)

for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    $firstCellText = $row.Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13)
    if ($firstCellText -eq "Status Cocktail") {
        $lastCell = $row.Cells.Item($row.Cells.Count)
        $cellRange = $lastCell.Range
        # Trim the trailing cell-mark / paragraph-mark characters before replacing.
        $cellRange.End = $cellRange.End - 1
        $cellRange.Text = "ok"
        break
    }
}
